# Edit: rename sheet, add a second ("deutsch") sheet, append daily-stats
# data for rows 126-174 on columns B:G, add a note in M164, and update
# the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "english"

# Add the new "deutsch" sheet right after "english".
# Two sheets are added in sequence (so the new one receives sheetId=3,
# matching the target workbook) and the first scratch sheet is then
# removed, leaving "deutsch" positioned immediately after "english".
$wsScratch = $wb.Worksheets.Add($null, $ws)
$wsDeutsch = $wb.Worksheets.Add($null, $wsScratch)
$wsDeutsch.Name = "deutsch"
$wsScratch.Delete()

# Fill in the daily statistics block (columns B:G) for rows 126-174.
$data = New-Object 'object[,]' 49,6
$data[0,0]=26; $data[0,1]=289; $data[0,2]=49; $data[0,3]=5505; $data[0,4]=41; $data[0,5]=125
$data[1,0]=2; $data[1,1]=243; $data[1,2]=35; $data[1,3]=5514; $data[1,4]=76; $data[1,5]=83
$data[2,0]=31; $data[2,1]=297; $data[2,2]=45; $data[2,3]=5540; $data[2,4]=163; $data[2,5]=$null
$data[3,0]=10; $data[3,1]=146; $data[3,2]=24; $data[3,3]=5542; $data[3,4]=84; $data[3,5]=78
$data[4,0]=11; $data[4,1]=230; $data[4,2]=30; $data[4,3]=5529; $data[4,4]=37; $data[4,5]=149
$data[5,0]=6; $data[5,1]=218; $data[5,2]=32; $data[5,3]=5565; $data[5,4]=80; $data[5,5]=76
$data[6,0]=4; $data[6,1]=210; $data[6,2]=33; $data[6,3]=5562; $data[6,4]=88; $data[6,5]=74
$data[7,0]=18; $data[7,1]=292; $data[7,2]=44; $data[7,3]=$null; $data[7,4]=$null; $data[7,5]=$null
$data[8,0]=12; $data[8,1]=228; $data[8,2]=36; $data[8,3]=5578; $data[8,4]=44; $data[8,5]=130
$data[9,0]=15; $data[9,1]=235; $data[9,2]=34; $data[9,3]=$null; $data[9,4]=$null; $data[9,5]=$null
$data[10,0]=11; $data[10,1]=212; $data[10,2]=32; $data[10,3]=5637; $data[10,4]=98; $data[10,5]=42
$data[11,0]=6; $data[11,1]=146; $data[11,2]=26; $data[11,3]=5691; $data[11,4]=74; $data[11,5]=18
$data[12,0]=1; $data[12,1]=157; $data[12,2]=22; $data[12,3]=5685; $data[12,4]=55; $data[12,5]=44
$data[13,0]=0; $data[13,1]=141; $data[13,2]=20; $data[13,3]=5673; $data[13,4]=75; $data[13,5]=36
$data[14,0]=3; $data[14,1]=132; $data[14,2]=20; $data[14,3]=5701; $data[14,4]=42; $data[14,5]=43
$data[15,0]=6; $data[15,1]=126; $data[15,2]=20; $data[15,3]=5697; $data[15,4]=44; $data[15,5]=51
$data[16,0]=7; $data[16,1]=130; $data[16,2]=20; $data[16,3]=5721; $data[16,4]=43; $data[16,5]=35
$data[17,0]=10; $data[17,1]=178; $data[17,2]=28; $data[17,3]=5719; $data[17,4]=55; $data[17,5]=35
$data[18,0]=0; $data[18,1]=138; $data[18,2]=21; $data[18,3]=5710; $data[18,4]=62; $data[18,5]=37
$data[19,0]=0; $data[19,1]=137; $data[19,2]=20; $data[19,3]=5741; $data[19,4]=49; $data[19,5]=19
$data[20,0]=1; $data[20,1]=72; $data[20,2]=10; $data[20,3]=5710; $data[20,4]=24; $data[20,5]=76
$data[21,0]=21; $data[21,1]=145; $data[21,2]=20; $data[21,3]=5714; $data[21,4]=47; $data[21,5]=70
$data[22,0]=1; $data[22,1]=116; $data[22,2]=22; $data[22,3]=5703; $data[22,4]=40; $data[22,5]=89
$data[23,0]=5; $data[23,1]=223; $data[23,2]=33; $data[23,3]=5724; $data[23,4]=36; $data[23,5]=77
$data[24,0]=11; $data[24,1]=145; $data[24,2]=23; $data[24,3]=5721; $data[24,4]=69; $data[24,5]=58
$data[25,0]=3; $data[25,1]=246; $data[25,2]=35; $data[25,3]=5678; $data[25,4]=58; $data[25,5]=115
$data[26,0]=1; $data[26,1]=168; $data[26,2]=24; $data[26,3]=5732; $data[26,4]=70; $data[26,5]=50
$data[27,0]=2; $data[27,1]=141; $data[27,2]=20; $data[27,3]=5715; $data[27,4]=71; $data[27,5]=68
$data[28,0]=12; $data[28,1]=142; $data[28,2]=18; $data[28,3]=5707; $data[28,4]=57; $data[28,5]=102
$data[29,0]=0; $data[29,1]=174; $data[29,2]=29; $data[29,3]=5767; $data[29,4]=68; $data[29,5]=31
$data[30,0]=2; $data[30,1]=150; $data[30,2]=20; $data[30,3]=5772; $data[30,4]=76; $data[30,5]=20
$data[31,0]=0; $data[31,1]=162; $data[31,2]=22; $data[31,3]=5746; $data[31,4]=72; $data[31,5]=59
$data[32,0]=8; $data[32,1]=134; $data[32,2]=20; $data[32,3]=5774; $data[32,4]=48; $data[32,5]=54
$data[33,0]=1; $data[33,1]=141; $data[33,2]=21; $data[33,3]=5781; $data[33,4]=48; $data[33,5]=48
$data[34,0]=1; $data[34,1]=137; $data[34,2]=20; $data[34,3]=5781; $data[34,4]=48; $data[34,5]=48
$data[35,0]=7; $data[35,1]=127; $data[35,2]=15; $data[35,3]=5802; $data[35,4]=49; $data[35,5]=33
$data[36,0]=2; $data[36,1]=110; $data[36,2]=16; $data[36,3]=5787; $data[36,4]=35; $data[36,5]=64
$data[37,0]=6; $data[37,1]=117; $data[37,2]=15; $data[37,3]=$null; $data[37,4]=$null; $data[37,5]=$null
$data[38,0]=24; $data[38,1]=180; $data[38,2]=38; $data[38,3]=$null; $data[38,4]=$null; $data[38,5]=$null
$data[39,0]=46; $data[39,1]=327; $data[39,2]=63; $data[39,3]=$null; $data[39,4]=$null; $data[39,5]=$null
$data[40,0]=11; $data[40,1]=216; $data[40,2]=43; $data[40,3]=$null; $data[40,4]=$null; $data[40,5]=$null
$data[41,0]=35; $data[41,1]=274; $data[41,2]=40; $data[41,3]=$null; $data[41,4]=$null; $data[41,5]=$null
$data[42,0]=2; $data[42,1]=190; $data[42,2]=31; $data[42,3]=$null; $data[42,4]=$null; $data[42,5]=$null
$data[43,0]=29; $data[43,1]=321; $data[43,2]=54; $data[43,3]=$null; $data[43,4]=$null; $data[43,5]=$null
$data[44,0]=18; $data[44,1]=314; $data[44,2]=51; $data[44,3]=$null; $data[44,4]=$null; $data[44,5]=$null
$data[45,0]=37; $data[45,1]=293; $data[45,2]=50; $data[45,3]=$null; $data[45,4]=$null; $data[45,5]=$null
$data[46,0]=10; $data[46,1]=270; $data[46,2]=47; $data[46,3]=$null; $data[46,4]=$null; $data[46,5]=$null
$data[47,0]=9; $data[47,1]=262; $data[47,2]=41; $data[47,3]=$null; $data[47,4]=$null; $data[47,5]=$null
$data[48,0]=3; $data[48,1]=150; $data[48,2]=36; $data[48,3]=$null; $data[48,4]=$null; $data[48,5]=$null
$ws.Range("B126:G174").Value = $data

# Note cell added in column M of row 164.
$ws.Range("M164").Value = "тут я почав додавати німецькі слова"

# Update the visible selection / scroll position on the "english" sheet.
$ws.Range("R170").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 145
$win.ScrollColumn = 1
